$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.500.16'
$ws.Range('E2').Value = '  +3.14%  '

$ws.Range('D3').Value = '1.604.32'
$ws.Range('E3').Value = '  +3.05%  '

$ws.Range('E4').Value = '  -0.03%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.07'
$ws.Range('E5').Value = '  +1.08%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.523'
$ws.Range('E6').Value = '  +7.71%  '

$ws.Range('E7').Value = '  -0.06%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '26.91'
$ws.Range('E8').Value = '  +9.63%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '43.36'
$ws.Range('E9').Value = '  -1.20%  '

$ws.Range('E10').Value = '  +2.07%  '

$ws.Range('E11').Value = '  +2.57%  '

$ws.Range('E12').Value = '  +1.90%  '

$ws.Range('D13').Value = '1.834.57'
$ws.Range('E13').Value = '  +3.04%  '

$ws.Range('D14').Value = '1.611.79'
$ws.Range('E14').Value = '  +3.40%  '

$ws.Range('D15').Value = '29.543.33'
$ws.Range('E15').Value = '  +3.17%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.536'
$ws.Range('E16').Value = '  +4.45%  '

$ws.Range('E17').Value = '  +2.73%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '63.51'
$ws.Range('E18').Value = '  +3.58%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '242.89'
$ws.Range('E19').Value = '  +5.60%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.61'

$ws.Range('D21').Value = '0.0₃0693'
$ws.Range('E21').Value = '  +3.14%  '

$ws.Range('E22').Value = '  +0.03%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.01'
$ws.Range('E23').Value = '  +2.58%  '

$ws.Range('E24').Value = '  +2.37%  '

$ws.Range('E25').Value = '  +0.47%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '154.39'
$ws.Range('E26').Value = '  +1.90%  '

$ws.Range('B27').Value = 'Stellar'
$ws.Range('C27').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.109'
$ws.Range('E27').Value = '  +6.18%  '

$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.30'
$ws.Range('E28').Value = '  +3.76%  '

$ws.Range('E29').Value = '  +2.68%  '

$ws.Range('E30').Value = '  +0.00%  '

$ws.Range('E31').Value = '  +2.89%  '

$ws.Range('E32').Value = '  -0.07%  '

$ws.Range('E33').Value = '  +1.86%  '

$ws.Range('E34').Value = '  +3.52%  '

$ws.Range('D35').Value = '1.416.84'
$ws.Range('E35').Value = '  +1.76%  '

$ws.Range('E36').Value = '  -2.35%  '

$ws.Range('E37').Value = '  +2.43%  '

$ws.Range('E38').Value = '  +5.36%  '

$ws.Range('E39').Value = '  +0.63%  '

$ws.Range('E40').Value = '  +2.01%  '

$ws.Range('E41').Value = '  +3.54%  '

$ws.Range('E42').Value = '  +1.27%  '

$ws.Range('B43').Value = 'PaxDollar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.998'
$ws.Range('E43').Value = '  -0.05%  '

$ws.Range('B44').Value = 'BitcoinSV'
$ws.Range('C44').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '52.86'
$ws.Range('E44').Value = '  +21.39%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.792'
$ws.Range('E45').Value = '  +2.22%  '

$ws.Range('E46').Value = '  +2.27%  '

$ws.Range('E47').Value = '  +2.62%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.28'
$ws.Range('E48').Value = '  -0.34%  '

$ws.Range('D49').Value = '1.745.52'
$ws.Range('E49').Value = '  +3.12%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '86.32'
$ws.Range('E50').Value = '  +1.30%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.835'
$ws.Range('E51').Value = '  -3.91%  '
